$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 322.25
$ws.Range("I5").Value = 339.7143
$ws.Range("K5").Value = 339.7143
$ws.Range("M5").Value = -224.7143
$ws.Range("H69").Value = 13196.637
$ws.Range("I69").Value = 15137.75
$ws.Range("J69").Value = 12087.429
$ws.Range("K69").Value = 45413.25
$ws.Range("L69").Value = 36262.287
$ws.Range("M69").Value = -44539.25
$ws.Range("N69").Value = -38010.287
$ws.Range("H72").Value = 13196.637
$ws.Range("I72").Value = 15137.75
$ws.Range("J72").Value = 12087.429
$ws.Range("K72").Value = 136239.75
$ws.Range("L72").Value = 108786.861
$ws.Range("M72").Value = -131871.75
$ws.Range("N72").Value = -117522.861
$ws.Range("H75").Value = 41248.75
$ws.Range("I75").Value = 29995
$ws.Range("K75").Value = 29995
$ws.Range("M75").Value = -29059
$ws.Range("H76").Value = 3527.111
$ws.Range("I76").Value = 3059
$ws.Range("K76").Value = 3059
$ws.Range("M76").Value = -2744
$ws.Range("H78").Value = 41248.75
$ws.Range("I78").Value = 29995
$ws.Range("K78").Value = 89985
$ws.Range("M78").Value = -85305
$ws.Range("H79").Value = 3527.111
$ws.Range("I79").Value = 3059
$ws.Range("K79").Value = 3059
$ws.Range("M79").Value = -1967
$ws.Range("H106").Value = 104370.3
$ws.Range("I106").Value = 147542.42
$ws.Range("J106").Value = 3635.3333
$ws.Range("K106").Value = 147542.42
$ws.Range("L106").Value = 3635.3333
$ws.Range("M106").Value = -146911.42
$ws.Range("N106").Value = -4897.3333
$ws.Range("H132").Value = 26319800
$ws.Range("I132").Value = 38465064
$ws.Range("J132").Value = 5056.25
$ws.Range("K132").Value = 115395192
$ws.Range("L132").Value = 15168.75
$ws.Range("M132").Value = -115392662
$ws.Range("N132").Value = -20228.75
$ws.Range("H137").Value = 1814.4073
$ws.Range("I137").Value = 1077.7646
$ws.Range("K137").Value = 3233.2938
$ws.Range("M137").Value = -683.2937999999999

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1939.7
$ws.Range("I2").Value = 1599.6666
$ws.Range("K2").Value = 1599.6666
$ws.Range("M2").Value = -1486.6666
$ws.Range("H80").Value = 59955.5
$ws.Range("J80").Value = 59955.5
$ws.Range("L80").Value = 59955.5
$ws.Range("N80").Value = -61951.5
$ws.Range("H83").Value = 59955.5
$ws.Range("J83").Value = 59955.5
$ws.Range("L83").Value = 179866.5
$ws.Range("N83").Value = -189850.5
$ws.Range("H116").Value = 1939.7
$ws.Range("I116").Value = 1599.6666
$ws.Range("K116").Value = 1599.6666
$ws.Range("M116").Value = 694.3334

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1939.7
$ws.Range("I3").Value = 1599.6666
$ws.Range("K3").Value = 1599.6666
$ws.Range("M3").Value = -1485.6666
$ws.Range("H35").Value = 3750
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H44").Value = 19166.666
$ws.Range("J44").Value = 19166.666
$ws.Range("L44").Value = 19166.666
$ws.Range("N44").Value = -20160.666
$ws.Range("H63").Value = 71107.2
$ws.Range("I63").Value = 35555
$ws.Range("J63").Value = 79995.25
$ws.Range("K63").Value = 35555
$ws.Range("L63").Value = 79995.25
$ws.Range("M63").Value = -34869
$ws.Range("N63").Value = -81367.25
$ws.Range("H66").Value = 71107.2
$ws.Range("I66").Value = 35555
$ws.Range("J66").Value = 79995.25
$ws.Range("K66").Value = 106665
$ws.Range("L66").Value = 239985.75
$ws.Range("M66").Value = -103233
$ws.Range("N66").Value = -246849.75
$ws.Range("H86").Value = 8335582.5
$ws.Range("I86").Value = 8335582.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 8335582.5
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -8334459.5
$ws.Range("H89").Value = 8335582.5
$ws.Range("I89").Value = 8335582.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 41677912.5
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -41672296.5
$ws.Range("H105").Value = 781.1429000000001
$ws.Range("I105").Value = 595.0769
$ws.Range("K105").Value = 595.0769
$ws.Range("M105").Value = 1151.9231
$ws.Range("H134").Value = 1569.3529
$ws.Range("I134").Value = 1360.7391
$ws.Range("K134").Value = 4082.2173
$ws.Range("M134").Value = -1547.2173

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H76").Value = 9998
$ws.Range("I76").Value = 9998
$ws.Range("K76").Value = 9998
$ws.Range("M76").Value = -9683
$ws.Range("H79").Value = 9998
$ws.Range("I79").Value = 9998
$ws.Range("K79").Value = 9998
$ws.Range("M79").Value = -8906
$ws.Range("H86").Value = 7916
$ws.Range("I86").Value = 8503.5
$ws.Range("J86").Value = 7622.25
$ws.Range("K86").Value = 8503.5
$ws.Range("L86").Value = 7622.25
$ws.Range("M86").Value = -7380.5
$ws.Range("N86").Value = -9868.25
$ws.Range("H88").Value = 55749.5
$ws.Range("J88").Value = 55749.5
$ws.Range("L88").Value = 55749.5
$ws.Range("N88").Value = -56561.5
$ws.Range("H89").Value = 7916
$ws.Range("I89").Value = 8503.5
$ws.Range("J89").Value = 7622.25
$ws.Range("K89").Value = 42517.5
$ws.Range("L89").Value = 38111.25
$ws.Range("M89").Value = -36901.5
$ws.Range("N89").Value = -49343.25
$ws.Range("H91").Value = 55749.5
$ws.Range("J91").Value = 55749.5
$ws.Range("L91").Value = 55749.5
$ws.Range("N91").Value = -58557.5
$ws.Range("H99").Value = 3988.923
$ws.Range("I99").Value = 3988.923
$ws.Range("K99").Value = 3988.923
$ws.Range("M99").Value = -2490.923
$ws.Range("H126").Value = 3988.923
$ws.Range("I126").Value = 3988.923
$ws.Range("K126").Value = 11966.769
$ws.Range("M126").Value = -9496.769
$ws.Range("H132").Value = 1918.5
$ws.Range("I132").Value = 1681.4615
$ws.Range("K132").Value = 5044.3845
$ws.Range("M132").Value = -2514.3845

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 134.73334
$ws.Range("J2").Value = 138
$ws.Range("L2").Value = 828
$ws.Range("N2").Value = -1054
$ws.Range("H44").Value = 738.5
$ws.Range("I44").Value = 880.55554
$ws.Range("J44").Value = 596.44446
$ws.Range("K44").Value = 2641.66662
$ws.Range("L44").Value = 1789.33338
$ws.Range("M44").Value = -2243.66662
$ws.Range("N44").Value = -2585.33338

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 1473.9231
$ws.Range("I13").Value = 325.83334
$ws.Range("J13").Value = 2458
$ws.Range("K13").Value = 325.83334
$ws.Range("L13").Value = 2458
$ws.Range("M13").Value = -186.83334
$ws.Range("N13").Value = -2736
$ws.Range("H57").Value = 51768.6
$ws.Range("J57").Value = 54940.31
$ws.Range("L57").Value = 54940.31
$ws.Range("N57").Value = -56580.31
$ws.Range("H70").Value = 14799
$ws.Range("I70").Value = 5998.75
$ws.Range("K70").Value = 5998.75
$ws.Range("M70").Value = -5728.75
$ws.Range("H73").Value = 14799
$ws.Range("I73").Value = 5998.75
$ws.Range("K73").Value = 5998.75
$ws.Range("M73").Value = -5062.75
$ws.Range("H102").Value = 63617.75
$ws.Range("I102").Value = 1038.8
$ws.Range("J102").Value = 167916
$ws.Range("K102").Value = 1038.8
$ws.Range("L102").Value = 167916
$ws.Range("M102").Value = 583.2
$ws.Range("N102").Value = -171160

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 5071.8887
$ws.Range("I32").Value = 5071.8887
$ws.Range("K32").Value = 5071.8887
$ws.Range("M32").Value = -4754.8887
$ws.Range("H98").Value = 74500
$ws.Range("J98").Value = 74500
$ws.Range("L98").Value = 74500
$ws.Range("N98").Value = -80490
$ws.Range("H132").Value = 5573.625
$ws.Range("I132").Value = 4397.75
$ws.Range("J132").Value = 6749.5
$ws.Range("K132").Value = 13193.25
$ws.Range("L132").Value = 20248.5
$ws.Range("M132").Value = -10663.25
$ws.Range("N132").Value = -25308.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4624.5815
$ws.Range("J122").Value = 4991.6
$ws.Range("L122").Value = 14974.8
$ws.Range("N122").Value = -19874.8
$ws.Range("H126").Value = 8983.186
$ws.Range("I126").Value = 7741.84
$ws.Range("J126").Value = 24500
$ws.Range("K126").Value = 23225.52
$ws.Range("L126").Value = 73500
$ws.Range("M126").Value = -20755.52
$ws.Range("N126").Value = -78440
$ws.Range("H132").Value = 1367.579
$ws.Range("I132").Value = 1163.9
$ws.Range("K132").Value = 3491.7
$ws.Range("M132").Value = -961.7000000000003
